$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants")

# xlLeft alignment constant
$xlLeft = -4131

# Row 19 - wbInit_Type (B keeps its existing left-aligned style, just update value)
$ws.Range("A19").Value = "wbInit_Type"
$ws.Range("B19").Value = "Task1, Initialization State"
$ws.Range("B19").HorizontalAlignment = $xlLeft
$ws.Range("C19").Value = "Name of Workblock"

# Row 20 - wbInit_SuppressSuccessful
$ws.Range("A20").Value = "wbInit_SuppressSuccessful"
$ws.Range("B20").Value = $false
$ws.Range("B20").HorizontalAlignment = $xlLeft
$ws.Range("C20").Value = "Do not log successful executions of wb"

# Row 21 - wbGetTransactionData_Type
$ws.Range("A21").Value = "wbGetTransactionData_Type"
$ws.Range("B21").Value = "Task1, Get Transaction Data State"
$ws.Range("B21").HorizontalAlignment = $xlLeft
$ws.Range("C21").Value = "Name of Workblock"

# Row 22 - wbGetTransactionData_SuppressSuccessful
$ws.Range("A22").Value = "wbGetTransactionData_SuppressSuccessful"
$ws.Range("B22").Value = $false
$ws.Range("B22").HorizontalAlignment = $xlLeft
$ws.Range("C22").Value = "Do not log successful executions of wb"

# Row 23 - wbProcessTransaction_Type
$ws.Range("A23").Value = "wbProcessTransaction_Type"
$ws.Range("B23").Value = "Task1, Process Transaction State"
$ws.Range("B23").HorizontalAlignment = $xlLeft
$ws.Range("C23").Value = "Name of Workblock"

# Row 24 - wbProcessTransaction_SuppressSuccessful
$ws.Range("A24").Value = "wbProcessTransaction_SuppressSuccessful"
$ws.Range("B24").Value = $false
$ws.Range("B24").HorizontalAlignment = $xlLeft
$ws.Range("C24").Value = "Do not log successful executions of wb"

# Row 25 - wbNextTransaction_Type
$ws.Range("A25").Value = "wbNextTransaction_Type"
$ws.Range("B25").Value = "Task1, Next Transaction"
$ws.Range("B25").HorizontalAlignment = $xlLeft
$ws.Range("C25").Value = "Name of Workblock"

# Row 26 - wbNextTransaction_SuppressSuccessful
$ws.Range("A26").Value = "wbNextTransaction_SuppressSuccessful"
$ws.Range("B26").Value = $false
$ws.Range("B26").HorizontalAlignment = $xlLeft
$ws.Range("C26").Value = "Do not log successful executions of wb"

# Row 27 - wbCloseAllApplications_Type
$ws.Range("A27").Value = "wbCloseAllApplications_Type"
$ws.Range("B27").Value = "Task1, Close All Applications"
$ws.Range("B27").HorizontalAlignment = $xlLeft
$ws.Range("C27").Value = "Name of Workblock"

# Row 28 - wbCloseAllApplications_SuppressSuccessful
$ws.Range("A28").Value = "wbCloseAllApplications_SuppressSuccessful"
$ws.Range("B28").Value = $false
$ws.Range("B28").HorizontalAlignment = $xlLeft
$ws.Range("C28").Value = "Do not log successful executions of wb"

# Row 29 - wbInitAllApplications_Type
$ws.Range("A29").Value = "wbInitAllApplications_Type"
$ws.Range("B29").Value = "Task1, InitAllApplications"
$ws.Range("B29").HorizontalAlignment = $xlLeft
$ws.Range("C29").Value = "Name of Workblock"

# Row 30 - wbInitAllApplications_SuppressSuccessful
$ws.Range("A30").Value = "wbInitAllApplications_SuppressSuccessful"
$ws.Range("B30").Value = $false
$ws.Range("B30").HorizontalAlignment = $xlLeft
$ws.Range("C30").Value = "Do not log successful executions of wb"

# Row 31 - wbProcess_Type
$ws.Range("A31").Value = "wbProcess_Type"
$ws.Range("B31").Value = "Task1, Process Transaction"
$ws.Range("B31").HorizontalAlignment = $xlLeft
$ws.Range("C31").Value = "Name of Workblock"

# Row 32 - wbProcess_SuppressSuccessful
$ws.Range("A32").Value = "wbProcess_SuppressSuccessful"
$ws.Range("B32").Value = $false
$ws.Range("B32").HorizontalAlignment = $xlLeft
$ws.Range("C32").Value = "Do not log successful executions of wb"

# Update selection to match the author's final cursor position
$ws.Activate()
$ws.Range("B26").Select()
